$d = $word.ActiveDocument

# Edit 1: "For football, the domain" -> "For soccer, the domain"
$d.Content.Find.Execute("For football, the domain", $true, $false, $false, $false, $false, $true, 1, $false, "For soccer, the domain", 2)

# Edit 2: "suitable model for football game prediction." -> "suitable model for soccer prediction."
$d.Content.Find.Execute("suitable model for football game prediction.", $true, $false, $false, $false, $false, $true, 1, $false, "suitable model for soccer prediction.", 2)
